$d = $word.ActiveDocument

# 1) Title line: "PSP 4.3 / Catering" -> "PSP 4.3 - Catering"
$found = $d.Content.Find.Execute("PSP 4.3 / Catering", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PSP 4.3 - Catering", 2)

# 2) Drop the stray "_GoBack" bookmark left over in the "Ergebnisse" paragraph
#    (text itself is unchanged - the bookmark just marked the author's last edit point)
If ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
